# Inform users with message boxes in case of errors or warnings
# - Adds new localization entries (AuthenticationFailed, ErrorDuringExecution,
#   ErrorDuringExecutionLog, CredentialNotFound, OperationCanceledByUser)
# - Tweaks wording for ConfirmNumerousRequests and GetSingleMachineFailure

$wb = $excel.ActiveWorkbook
$wsSettings = $wb.Worksheets.Item("Settings")
$wsAdvanced = $wb.Worksheets.Item("Advanced Settings")
$wsLocalization = $wb.Worksheets.Item("Localization")

$lo = $wsLocalization.ListObjects.Item("Table13")

# 1) Insert 3 new rows right after "TokenNotRetrieved" (before old row 63 "TokenExpired")
$wsLocalization.Rows("63:65").Insert()

$wsLocalization.Range("A63").Value = "AuthenticationFailed"
$wsLocalization.Range("B63").Value = "Authentication failed. Please check logs for more details."
$wsLocalization.Range("C63").Value = "認証に失敗しました。詳細はログを確認してください。"

$wsLocalization.Range("A64").Value = "ErrorDuringExecutionLog"
$wsLocalization.Range("B64").Value = "There was an execution error: {0} at {1}."
$wsLocalization.Range("C64").Value = "実行中にエラーが発生しました：{1}で{0}。"

$wsLocalization.Range("A65").Value = "ErrorDuringExecution"
$wsLocalization.Range("B65").Value = "There was an execution error. Please check logs for more details."
$wsLocalization.Range("C65").Value = "実行中にエラーが発生しました。詳細はログを確認してください。"

# 2) Insert 2 new rows right after "TypeNotSupported" (now at row 80, before the blank row)
$wsLocalization.Rows("81:82").Insert()

$wsLocalization.Range("A81").Value = "CredentialNotFound"
$wsLocalization.Range("B81").Value = "The specified credential was not found."
$wsLocalization.Range("C81").Value = "指定された資格情報が見つかりませんでした。"

$wsLocalization.Range("A82").Value = "OperationCanceledByUser"
$wsLocalization.Range("B82").Value = "Operation canceled by user."
$wsLocalization.Range("C82").Value = "ユーザーが操作をキャンセルしました。"

# 3) Update wording of the existing ConfirmNumerousRequests message (now at row 71)
$wsLocalization.Range("B71").Value = "The selected operation will make a large number of HTTP requests that might impact on Orchestrator's infrastructure. Continue the processing?"

# 4) Update wording of the existing GetSingleMachineFailure message (now at row 113): "Id" -> "ID"
$wsLocalization.Range("B113").Value = "Failed to get machine with ID: {0}. Request status: {1} / Response: {2}."

# 5) Grow the Table13 range/autofilter to cover the 5 newly inserted rows
$lo.Resize($wsLocalization.Range("A1:C156"))

# 6) Normalize the duplicated selection reference on the Advanced Settings / Localization sheets
$wsAdvanced.Activate()
$wsAdvanced.Range("A2").Select()

$wsLocalization.Activate()
$wsLocalization.Range("A2").Select()

# 7) Restore the originally active sheet
$wsSettings.Activate()
